$wb = $excel.ActiveWorkbook

# --- Sheet 1 updates ---
$ws1 = $wb.Worksheets.Item("Sheet 1")

# Force these cells to remain plain text so Excel does not
# auto-convert date/time-looking strings into date/time serials.
$ws1.Range("B2:B4").NumberFormat = "@"
$ws1.Range("B8").NumberFormat = "@"

$ws1.Range("B2").Value = "Thanush_Test"
$ws1.Range("B3").Value = "10/05/1999"
$ws1.Range("B4").Value = "03:04:00"
$ws1.Range("B8").Value = "Asia/Kolkata"

# --- Sheet 2 updates ---
$ws2 = $wb.Worksheets.Item("Sheet 2")

# Header
$ws2.Range("G2").Value = "Retro(R)"

# Row 3 - Ascendant
$ws2.Range("B3").Value = "Pisces"
$ws2.Range("C3").Value = "Jupiter"
$ws2.Range("D3").Value = "Purva Bhadrapada"
$ws2.Range("E3").Value = "Jupiter"
$ws2.Range("F3").Value = 1.263900747942614

# Row 4 - Sun
$ws2.Range("B4").Value = "Aries"
$ws2.Range("C4").Value = "Mars"
$ws2.Range("D4").Value = "Bharani"
$ws2.Range("E4").Value = "Venus"
$ws2.Range("F4").Value = 24.96816340468312
$ws2.Range("J4").Value = 2

# Row 5 - Moon
$ws2.Range("B5").Value = "Aquarius"
$ws2.Range("C5").Value = "Saturn"
$ws2.Range("D5").Value = "Shatabhisha"
$ws2.Range("E5").Value = "Rahu"
$ws2.Range("F5").Value = 309.075052640371
$ws2.Range("J5").Value = 12

# Row 6 - Mercury
$ws2.Range("B6").Value = "Aries"
$ws2.Range("C6").Value = "Mars"
$ws2.Range("D6").Value = "Ashwini"
$ws2.Range("E6").Value = "Ketu"
$ws2.Range("F6").Value = 7.96109256468479

# Row 7 - Venus
$ws2.Range("B7").Value = "Gemini"
$ws2.Range("C7").Value = "Mercury"
$ws2.Range("D7").Value = "Ardra"
$ws2.Range("E7").Value = "Rahu"
$ws2.Range("F7").Value = 67.50609350739118
$ws2.Range("G7").Value = "Direct"
$ws2.Range("H7").Value = "No"
$ws2.Range("J7").Value = 4

# Row 8 - Mars
$ws2.Range("D8").Value = "Chitra"
$ws2.Range("E8").Value = "Mars"
$ws2.Range("F8").Value = 184.8265144931989
$ws2.Range("G8").Value = "Retro"
$ws2.Range("J8").Value = 8

# Row 9 - Jupiter
$ws2.Range("B9").Value = "Pisces"
$ws2.Range("C9").Value = "Jupiter"
$ws2.Range("D9").Value = "Revati"
$ws2.Range("E9").Value = "Mercury"
$ws2.Range("F9").Value = 356.3771724279898
$ws2.Range("J9").Value = 1

# Row 10 - Saturn
$ws2.Range("F10").Value = 14.49292394624836
$ws2.Range("H10").Value = "Combust"
$ws2.Range("J10").Value = 2

# Row 11 - Uranus
$ws2.Range("F11").Value = 292.8890339183791
$ws2.Range("G11").Value = "Direct"
$ws2.Range("J11").Value = 11

# Row 12 - Neptune
$ws2.Range("D12").Value = "Shravana"
$ws2.Range("E12").Value = "Moon"
$ws2.Range("F12").Value = 280.5216800641753
$ws2.Range("J12").Value = 11

# Row 13 - Pluto
$ws2.Range("F13").Value = 225.8449774019211
$ws2.Range("J13").Value = 9

# Row 14 - Rahu
$ws2.Range("F14").Value = 113.7254071935777
$ws2.Range("J14").Value = 5

# Row 15 - Ketu
$ws2.Range("D15").Value = "Dhanishta"
$ws2.Range("E15").Value = "Mars"
$ws2.Range("F15").Value = 293.560692356362
$ws2.Range("J15").Value = 11
